$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = "ejr;lkwjlkw``kjkejg"
$ws.Range("J2").Value = "kj;lknds"
$ws.Range("K2").Value = "n;lkngdlkndlkgn;lkdg"

$ws.Range("K2").Select()
